$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: GPS Receiver (first table) ---
$ws.Range("A4").Value = "GPS Receiver"
$ws.Range("B4").Value = "OEM 615"
$ws.Range("C4").Value = "Novatel"
$ws.Range("D4").Value = "1W"
$ws.Range("E4").Value = "N/A"
$ws.Range("F4").Value = "Read GPS data and calculate Positon"
$ws.Range("G4").Value = "PDF link"

# --- Row 18: GPS Receiver (second table, duplicate of row 4) ---
$ws.Range("A18").Value = "GPS Receiver"
$ws.Range("B18").Value = "OEM 615"
$ws.Range("C18").Value = "Novatel"
$ws.Range("D18").Value = "1W"
$ws.Range("E18").Value = "N/A"
$ws.Range("F18").Value = "Read GPS data and calculate Positon"
$ws.Range("G18").Value = "PDF link"

# --- Row 19: GPS Antenna (second table - Antdevco option) ---
$ws.Range("A19").Value = "GPS Antenna"
$ws.Range("B19").Value = "Single Frequency Mircostrip"
$ws.Range("C19").Value = "Antdevco"
$ws.Range("D19").Value = "N/A"
$ws.Range("E19").Value = 5000
$ws.Range("F19").Value = "Recieves GPS and passes Signal to GPS"
$ws.Range("G19").Value = "PDF link"

# --- Row 5: GPS Antenna (first table - SpaceQuest option) ---
$ws.Range("A5").Value = "GPS Antenna"
$ws.Range("B5").Value = "ANT-GPS"
$ws.Range("C5").Value = "SpaceQuest"
$ws.Range("D5").Value = "N/A"
$ws.Range("E5").Value = 3000
$ws.Range("F5").Value = "Recieves GPS and passes Signal to GPS"
$ws.Range("G5").Value = "PDF link"

# --- Hyperlinks: add in G4, G18, G19, G5 order (matches rId3..rId6) ---
$ws.Hyperlinks.Add($ws.Range("G4"), "http://www.novatel.com/assets/documents/papers/OEM615.pdf")
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("G18"), "http://www.novatel.com/assets/documents/papers/OEM615.pdf")
$ws.Range("G3").Copy()
$ws.Range("G18").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("G19"), "http://www.antdevco.com/pdf/ant-gps.pdf")
$ws.Range("G3").Copy()
$ws.Range("G19").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("G5"), "http://www.spacequest.com/antennas/ANT-GPS.pdf")
$ws.Range("G3").Copy()
$ws.Range("G5").PasteSpecial(-4122)

# --- Selection / view state ---
$ws.Range("G5").Select()
